# Añadiendo Serie2_PEA_MX: completa las filas 68-73 (fechas trimestrales +
# valores de PEA) que estaban vacías (solo conservaban el formato de celda).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Datos nuevos: fecha (columna A, formato fecha ya existente) y PEA (columna B).
$newData = @(
    @{ Row = 68; Date = "07/01/2021"; Value = 58307446 },
    @{ Row = 69; Date = "10/01/2021"; Value = 58761793 },
    @{ Row = 70; Date = "01/01/2022"; Value = 58085314 },
    @{ Row = 71; Date = "04/01/2022"; Value = 59338419 },
    @{ Row = 72; Date = "07/01/2022"; Value = 59480471 },
    @{ Row = 73; Date = "10/01/2022"; Value = 60145456 }
)

foreach ($item in $newData) {
    $r = $item.Row

    # Columna A conserva el formato de fecha (estilo ya aplicado en la celda).
    $ws.Cells.Item($r, 1).Value = $item.Date

    # Columna B pasa a formato general (sin separador de miles) al rellenarla.
    $cellB = $ws.Cells.Item($r, 2)
    $cellB.Style = "Normal"
    $cellB.Value = $item.Value
}
